# Add two new worksheets ("ODI Batting Extra" and "ODI Bowling Extra") to the
# Afghanistan PlayerPerformance workbook, mirroring the header style of the
# existing sheets and filling in the single data row for match 4525.

$wb = $excel.ActiveWorkbook

# Reference sheet whose header formatting (bold font, borders, centered
# alignment) we want to copy onto the header rows of the new sheets.
$refSheet = $wb.Worksheets.Item(1)

# --- "ODI Batting Extra" ---------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$battingExtra = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$battingExtra.Name = "ODI Batting Extra"

$refSheet.Range("A1:D1").Copy()
$battingExtra.Range("A1:F1").PasteSpecial(-4122)

$battingExtra.Range("A1").Value = "MATCH_CODE"
$battingExtra.Range("B1").Value = "BATTING_POSITION"
$battingExtra.Range("C1").Value = "NUM_4"
$battingExtra.Range("D1").Value = "NUM_6"
$battingExtra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$battingExtra.Range("F1").Value = "MAN_OF_MATCH"

$battingExtra.Range("A2").Value = "4525"
$battingExtra.Range("F2").Value = "NO"

# --- "ODI Bowling Extra" ----------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$bowlingExtra = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$bowlingExtra.Name = "ODI Bowling Extra"

$refSheet.Range("A1:C1").Copy()
$bowlingExtra.Range("A1:C1").PasteSpecial(-4122)

$bowlingExtra.Range("A1").Value = "MATCH_CODE"
$bowlingExtra.Range("B1").Value = "MAIDEN_OVERS"
$bowlingExtra.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

$bowlingExtra.Range("A2").Value = "4525"

$refSheet.Select()
